$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates ---
$ws.Range("A8").Value = "Volume 32   Number  14"
$ws.Range("C9").Value = "Report Covering the Week  3/31/2025  Through  4/6/2025"

# --- Column E width adjustment ---
$ws.Columns("E").ColumnWidth = 6.168446

# --- Cells changing type (numeric <-> text placeholder) : copy style+type from template cells, then fix numeric value where needed ---
$ws.Range("C14").Copy($ws.Range("C15"))
$ws.Range("C14").Copy($ws.Range("G16"))
$ws.Range("E14").Copy($ws.Range("H16"))
$ws.Range("C14").Copy($ws.Range("D18"))
$ws.Range("E14").Copy($ws.Range("E18"))
$ws.Range("J14").Copy($ws.Range("C20"))
$ws.Range("C20").Value = 1
$ws.Range("J14").Copy($ws.Range("D20"))
$ws.Range("D20").Value = 2
$ws.Range("K14").Copy($ws.Range("E20"))
$ws.Range("E20").Value = -50
$ws.Range("C14").Copy($ws.Range("C27"))
$ws.Range("C14").Copy($ws.Range("D28"))
$ws.Range("E14").Copy($ws.Range("E28"))
$ws.Range("C14").Copy($ws.Range("D31"))
$ws.Range("E14").Copy($ws.Range("E31"))

# --- Simple numeric value updates (style/type unchanged) ---
$ws.Range("N16").Value = -73.333333333333
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 4
$ws.Range("E17").Value = -25
$ws.Range("I17").Value = 31
$ws.Range("J17").Value = 35
$ws.Range("K17").Value = -11.428571428571
$ws.Range("L17").Value = 55
$ws.Range("M17").Value = 158.333333333333
$ws.Range("N17").Value = 72.222222222222
$ws.Range("I18").Value = 20
$ws.Range("K18").Value = 300
$ws.Range("L18").Value = 5.263157894736
$ws.Range("M18").Value = -28.571428571428
$ws.Range("N18").Value = -69.230769230769
$ws.Range("D19").Value = 4
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 21
$ws.Range("G19").Value = 20
$ws.Range("H19").Value = 5
$ws.Range("I19").Value = 61
$ws.Range("J19").Value = 72
$ws.Range("K19").Value = -15.277777777777
$ws.Range("L19").Value = -17.567567567567
$ws.Range("M19").Value = 45.238095238095
$ws.Range("N19").Value = 69.444444444444
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 3
$ws.Range("H20").Value = -33.333333333333
$ws.Range("I20").Value = 7
$ws.Range("J20").Value = 13
$ws.Range("K20").Value = -46.153846153846
$ws.Range("L20").Value = -63.157894736842
$ws.Range("M20").Value = -30
$ws.Range("N20").Value = -96.067415730337
$ws.Range("C21").Value = 10
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 0
$ws.Range("F21").Value = 35
$ws.Range("H21").Value = -5.405405405405
$ws.Range("I21").Value = 125
$ws.Range("J21").Value = 134
$ws.Range("K21").Value = -6.716417910447
$ws.Range("L21").Value = -8.759124087591
$ws.Range("M21").Value = 31.578947368421
$ws.Range("N21").Value = -60.31746031746
$ws.Range("C24").Value = 7
$ws.Range("E24").Value = 16.666666666666
$ws.Range("F24").Value = 25
$ws.Range("G24").Value = 25
$ws.Range("H24").Value = 0
$ws.Range("I24").Value = 105
$ws.Range("J24").Value = 105
$ws.Range("K24").Value = 0
$ws.Range("L24").Value = -26.573426573426
$ws.Range("M24").Value = -9.482758620689
$ws.Range("C25").Value = 2
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -33.333333333333
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -36.842105263157
$ws.Range("I25").Value = 52
$ws.Range("J25").Value = 53
$ws.Range("K25").Value = -1.88679245283
$ws.Range("L25").Value = -36.585365853658
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 15
$ws.Range("G26").Value = 11
$ws.Range("H26").Value = 36.363636363636
$ws.Range("I26").Value = 61
$ws.Range("J26").Value = 36
$ws.Range("K26").Value = 69.444444444444
$ws.Range("L26").Value = 17.307692307692
$ws.Range("M26").Value = 19.607843137254
